{"js": "// Replace each three-digit-divided-by-one-digit expression with its updated value.\n// Mapping derived from the authoritative OOXML diff (old text -> new text),\n// applied in document order; every old string is unique in this document.\nconst replacements = [\n  [\"875\u00f77=\", \"120\u00f76=\"],\n  [\"244\u00f72=\", \"122\u00f76=\"],\n  [\"740\u00f72=\", \"856\u00f79=\"],\n  [\"973\u00f74=\", \"383\u00f75=\"],\n  [\"222\u00f79=\", \"272\u00f75=\"],\n  [\"838\u00f75=\", \"876\u00f75=\"],\n  [\"770\u00f73=\", \"225\u00f75=\"],\n  [\"221\u00f77=\", \"701\u00f74=\"],\n  [\"874\u00f78=\", \"251\u00f72=\"],\n  [\"982\u00f78=\", \"679\u00f73=\"],\n  [\"387\u00f72=\", \"189\u00f75=\"],\n  [\"597\u00f78=\", \"945\u00f75=\"],\n  [\"904\u00f78=\", \"870\u00f75=\"],\n  [\"163\u00f76=\", \"754\u00f77=\"],\n  [\"131\u00f79=\", \"182\u00f75=\"],\n  [\"176\u00f78=\", \"675\u00f76=\"],\n  [\"988\u00f74=\", \"344\u00f76=\"],\n  [\"346\u00f72=\", \"881\u00f72=\"],\n  [\"509\u00f75=\", \"208\u00f77=\"],\n  [\"878\u00f75=\", \"418\u00f73=\"],\n  [\"994\u00f75=\", \"101\u00f72=\"],\n  [\"166\u00f76=\", \"846\u00f74=\"],\n  [\"138\u00f78=\", \"462\u00f78=\"],\n  [\"137\u00f74=\", \"546\u00f74=\"],\n  [\"690\u00f79=\", \"655\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-divided-by-one-digit expression with its updated value.\n# Mapping derived from the authoritative OOXML diff (old text -> new text).\n# Every old string is unique in this document, so Find/Replace is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"875\u00f77=\", \"120\u00f76=\"),\n    @(\"244\u00f72=\", \"122\u00f76=\"),\n    @(\"740\u00f72=\", \"856\u00f79=\"),\n    @(\"973\u00f74=\", \"383\u00f75=\"),\n    @(\"222\u00f79=\", \"272\u00f75=\"),\n    @(\"838\u00f75=\", \"876\u00f75=\"),\n    @(\"770\u00f73=\", \"225\u00f75=\"),\n    @(\"221\u00f77=\", \"701\u00f74=\"),\n    @(\"874\u00f78=\", \"251\u00f72=\"),\n    @(\"982\u00f78=\", \"679\u00f73=\"),\n    @(\"387\u00f72=\", \"189\u00f75=\"),\n    @(\"597\u00f78=\", \"945\u00f75=\"),\n    @(\"904\u00f78=\", \"870\u00f75=\"),\n    @(\"163\u00f76=\", \"754\u00f77=\"),\n    @(\"131\u00f79=\", \"182\u00f75=\"),\n    @(\"176\u00f78=\", \"675\u00f76=\"),\n    @(\"988\u00f74=\", \"344\u00f76=\"),\n    @(\"346\u00f72=\", \"881\u00f72=\"),\n    @(\"509\u00f75=\", \"208\u00f77=\"),\n    @(\"878\u00f75=\", \"418\u00f73=\"),\n    @(\"994\u00f75=\", \"101\u00f72=\"),\n    @(\"166\u00f76=\", \"846\u00f74=\"),\n    @(\"138\u00f78=\", \"462\u00f78=\"),\n    @(\"137\u00f74=\", \"546\u00f74=\"),\n    @(\"690\u00f79=\", \"655\u00f74=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
